# Auto-generated Excel COM-interop script to apply the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 8065781  # H33: was 8622031
$ws.Cells.Item(33, 9).Value = 11363899  # I33: was 12500272
$ws.Cells.Item(33, 11).Value = 11363899  # K33: was 12500272
$ws.Cells.Item(33, 13).Value = -11363670  # M33: was -12500043

$ws.Cells.Item(112, 8).Value = 5328.879  # H112: was 5565.9355
$ws.Cells.Item(112, 10).Value = 9972.235000000001  # J112: was 11081.267
$ws.Cells.Item(112, 12).Value = 29916.705  # L112: was 33243.801
$ws.Cells.Item(112, 14).Value = -32132.705  # N112: was -35459.801

$ws.Cells.Item(113, 8).Value = 8619.786  # H113: was 9205.923000000001
$ws.Cells.Item(113, 10).Value = 4250  # J113: was 7500
$ws.Cells.Item(113, 12).Value = 4250  # L113: was 7500
$ws.Cells.Item(113, 14).Value = -10758  # N113: was -14008

$ws.Cells.Item(116, 8).Value = 5930891  # H116: was 7412860
$ws.Cells.Item(116, 9).Value = 8893988  # I116: was 11116729
$ws.Cells.Item(116, 10).Value = 4698  # J116: was 5122.5
$ws.Cells.Item(116, 11).Value = 8893988  # K116: was 11116729
$ws.Cells.Item(116, 12).Value = 4698  # L116: was 5122.5
$ws.Cells.Item(116, 13).Value = -8890546  # M116: was -11113287
$ws.Cells.Item(116, 14).Value = -11582  # N116: was -12006.5

$ws.Cells.Item(132, 8).Value = 12964.028  # H132: was 13708.667
$ws.Cells.Item(132, 9).Value = 7385.9395  # I132: was 7818.7417
$ws.Cells.Item(132, 11).Value = 22157.8185  # K132: was 23456.2251
$ws.Cells.Item(132, 13).Value = -19627.8185  # M132: was -20926.2251

$ws.Cells.Item(137, 8).Value = 202525.64  # H137: was 222659.4
$ws.Cells.Item(137, 9).Value = 3472.125  # I137: was 3931.5
$ws.Cells.Item(137, 10).Value = 733335  # J137: was 550751.25
$ws.Cells.Item(137, 11).Value = 10416.375  # K137: was 11794.5
$ws.Cells.Item(137, 12).Value = 2200005  # L137: was 1652253.75
$ws.Cells.Item(137, 13).Value = -7866.375  # M137: was -9244.5
$ws.Cells.Item(137, 14).Value = -2205105  # N137: was -1657353.75

$ws.Cells.Item(138, 8).Value = 3789.963  # H138: was 3847.2886
$ws.Cells.Item(138, 9).Value = 3431.5  # I138: was 3452.7856
$ws.Cells.Item(138, 10).Value = 3915.425  # J138: was 3992.6316
$ws.Cells.Item(138, 11).Value = 10294.5  # K138: was 10358.3568
$ws.Cells.Item(138, 12).Value = 11746.275  # L138: was 11977.8948
$ws.Cells.Item(138, 13).Value = -5154.5  # M138: was -5218.356800000001
$ws.Cells.Item(138, 14).Value = -22026.275  # N138: was -22257.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 23174.053  # H74: was 24433.445
$ws.Cells.Item(74, 9).Value = 1356.6666  # I74: was 1427.5
$ws.Cells.Item(74, 10).Value = 60575.285  # J74: was 53190.875
$ws.Cells.Item(74, 11).Value = 1356.6666  # K74: was 1427.5
$ws.Cells.Item(74, 12).Value = 60575.285  # L74: was 53190.875
$ws.Cells.Item(74, 13).Value = -482.6666  # M74: was -553.5
$ws.Cells.Item(74, 14).Value = -62323.285  # N74: was -54938.875

$ws.Cells.Item(77, 8).Value = 23174.053  # H77: was 24433.445
$ws.Cells.Item(77, 9).Value = 1356.6666  # I77: was 1427.5
$ws.Cells.Item(77, 10).Value = 60575.285  # J77: was 53190.875
$ws.Cells.Item(77, 11).Value = 6783.333000000001  # K77: was 7137.5
$ws.Cells.Item(77, 12).Value = 302876.425  # L77: was 265954.375
$ws.Cells.Item(77, 13).Value = -2415.333000000001  # M77: was -2769.5
$ws.Cells.Item(77, 14).Value = -311612.425  # N77: was -274690.375

$ws.Cells.Item(102, 8).Value = 4466.3335  # H102: was 3624.5
$ws.Cells.Item(102, 9).Value = 5002.25  # I102: was 3286.7144
$ws.Cells.Item(102, 10).Value = 3394.5  # J102: was 5989
$ws.Cells.Item(102, 11).Value = 5002.25  # K102: was 3286.7144
$ws.Cells.Item(102, 12).Value = 3394.5  # L102: was 5989
$ws.Cells.Item(102, 13).Value = -3380.25  # M102: was -1664.7144
$ws.Cells.Item(102, 14).Value = -6638.5  # N102: was -9233

$ws.Cells.Item(110, 8).Value = 7281.875  # H110: was 8222.143
$ws.Cells.Item(110, 9).Value = 9117  # I110: was 10800.4
$ws.Cells.Item(110, 11).Value = 9117  # K110: was 10800.4
$ws.Cells.Item(110, 13).Value = -7072  # M110: was -8755.4

$ws.Cells.Item(122, 8).Value = 1142.6666  # H122: was 1175.1177
$ws.Cells.Item(122, 9).Value = 1023.625  # I122: was 1052.4667
$ws.Cells.Item(122, 11).Value = 3070.875  # K122: was 3157.4001
$ws.Cells.Item(122, 13).Value = -620.875  # M122: was -707.4000999999998

$ws.Cells.Item(132, 8).Value = 3313062  # H132: was 3581668
$ws.Cells.Item(132, 9).Value = 1655.5758  # I132: was 1814.6207
$ws.Cells.Item(132, 10).Value = 18923978  # J132: was 16558636
$ws.Cells.Item(132, 11).Value = 4966.7274  # K132: was 5443.8621
$ws.Cells.Item(132, 12).Value = 56771934  # L132: was 49675908
$ws.Cells.Item(132, 13).Value = -2436.7274  # M132: was -2913.8621
$ws.Cells.Item(132, 14).Value = -56776994  # N132: was -49680968

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1762.7142  # H22: was 519.7143
$ws.Cells.Item(22, 9).Value = 2610  # I22: was 434.75
$ws.Cells.Item(22, 11).Value = 2610  # K22: was 434.75
$ws.Cells.Item(22, 13).Value = -2437  # M22: was -261.75

$ws.Cells.Item(107, 8).Value = 546.3333  # H107: was 379.75
$ws.Cells.Item(107, 9).Value = 546.3333  # I107: was 362.64285
$ws.Cells.Item(107, 10).Value = 0  # J107: was 499.5
$ws.Cells.Item(107, 11).Value = 546.3333  # K107: was 362.64285
$ws.Cells.Item(107, 12).Value = 0  # L107: was 499.5
$ws.Cells.Item(107, 13).Value = 1373.6667  # M107: was 1557.35715
$ws.Cells.Item(107, 14).ClearContents()  # N107: was -4339.5

$ws.Cells.Item(134, 8).Value = 63305  # H134: was 86592.44
$ws.Cells.Item(134, 9).Value = 78120.14  # I134: was 121272.336
$ws.Cells.Item(134, 10).Value = 37378.5  # J134: was 42004
$ws.Cells.Item(134, 11).Value = 234360.42  # K134: was 363817.008
$ws.Cells.Item(134, 12).Value = 112135.5  # L134: was 126012
$ws.Cells.Item(134, 13).Value = -231825.42  # M134: was -361282.008
$ws.Cells.Item(134, 14).Value = -117205.5  # N134: was -131082

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 130582.94  # H31: was 138669.69
$ws.Cells.Item(31, 9).Value = 134658.8  # I31: was 144191.92
$ws.Cells.Item(31, 11).Value = 134658.8  # K31: was 144191.92
$ws.Cells.Item(31, 13).Value = -134363.8  # M31: was -143896.92

$ws.Cells.Item(34, 8).Value = 130582.94  # H34: was 138669.69
$ws.Cells.Item(34, 9).Value = 134658.8  # I34: was 144191.92
$ws.Cells.Item(34, 11).Value = 134658.8  # K34: was 144191.92
$ws.Cells.Item(34, 13).Value = -134456.8  # M34: was -143989.92

$ws.Cells.Item(58, 8).Value = 12166.828  # H58: was 9662.621999999999
$ws.Cells.Item(58, 9).Value = 4357.2856  # I58: was 3446.9211
$ws.Cells.Item(58, 11).Value = 4357.2856  # K58: was 3446.9211
$ws.Cells.Item(58, 13).Value = -4154.2856  # M58: was -3243.9211

$ws.Cells.Item(132, 8).Value = 85960180  # H132: was 103152070
$ws.Cells.Item(132, 9).Value = 233152.33  # I132: was 299568.84
$ws.Cells.Item(132, 11).Value = 699456.99  # K132: was 898706.52
$ws.Cells.Item(132, 13).Value = -696926.99  # M132: was -896176.52

$ws.Cells.Item(136, 8).Value = 12166.828  # H136: was 9662.621999999999
$ws.Cells.Item(136, 9).Value = 4357.2856  # I136: was 3446.9211
$ws.Cells.Item(136, 11).Value = 13071.8568  # K136: was 10340.7633
$ws.Cells.Item(136, 13).Value = -10521.8568  # M136: was -7790.763300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1200  # H107: was 1134.5714
$ws.Cells.Item(107, 9).Value = 1300  # I107: was 1026.6666
$ws.Cells.Item(107, 10).Value = 1190  # J107: was 1164
$ws.Cells.Item(107, 11).Value = 3900  # K107: was 3079.9998
$ws.Cells.Item(107, 12).Value = 3570  # L107: was 3492
$ws.Cells.Item(107, 13).Value = -1980  # M107: was -1159.9998
$ws.Cells.Item(107, 14).Value = -7410  # N107: was -7332

$ws.Cells.Item(132, 8).Value = 1784145.6  # H132: was 2166044.5
$ws.Cells.Item(132, 10).Value = 5051960  # J132: was 10101969
$ws.Cells.Item(132, 12).Value = 45467640  # L132: was 90917721
$ws.Cells.Item(132, 14).Value = -45472700  # N132: was -90922781

$ws.Cells.Item(137, 8).Value = 8483.5  # H137: was 7806.1816
$ws.Cells.Item(137, 10).Value = 11257.857  # J137: was 9979.75
$ws.Cells.Item(137, 12).Value = 33773.571  # L137: was 29939.25
$ws.Cells.Item(137, 14).Value = -43973.571  # N137: was -40139.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2913.7144  # H80: was 3205
$ws.Cells.Item(80, 9).Value = 2877.7856  # I80: was 3022.2307
$ws.Cells.Item(80, 10).Value = 2985.5715  # J80: was 3680.2
$ws.Cells.Item(80, 11).Value = 2877.7856  # K80: was 3022.2307
$ws.Cells.Item(80, 12).Value = 2985.5715  # L80: was 3680.2
$ws.Cells.Item(80, 13).Value = -1879.7856  # M80: was -2024.2307
$ws.Cells.Item(80, 14).Value = -4981.5715  # N80: was -5676.2

$ws.Cells.Item(83, 8).Value = 2913.7144  # H83: was 3205
$ws.Cells.Item(83, 9).Value = 2877.7856  # I83: was 3022.2307
$ws.Cells.Item(83, 10).Value = 2985.5715  # J83: was 3680.2
$ws.Cells.Item(83, 11).Value = 14388.928  # K83: was 15111.1535
$ws.Cells.Item(83, 12).Value = 14927.8575  # L83: was 18401
$ws.Cells.Item(83, 13).Value = -9396.928  # M83: was -10119.1535
$ws.Cells.Item(83, 14).Value = -24911.8575  # N83: was -28385

$ws.Cells.Item(122, 8).Value = 2726.1304  # H122: was 2997.0557
$ws.Cells.Item(122, 9).Value = 2464  # I122: was 2726.7693
$ws.Cells.Item(122, 10).Value = 3669.8  # J122: was 3699.8
$ws.Cells.Item(122, 11).Value = 7392  # K122: was 8180.3079
$ws.Cells.Item(122, 12).Value = 11009.4  # L122: was 11099.4
$ws.Cells.Item(122, 13).Value = -4942  # M122: was -5730.3079
$ws.Cells.Item(122, 14).Value = -15909.4  # N122: was -15999.4

$ws.Cells.Item(132, 8).Value = 949457.4399999999  # H132: was 1068089.8
$ws.Cells.Item(132, 9).Value = 933.6  # I132: was 1035.875
$ws.Cells.Item(132, 10).Value = 2135112.2  # J132: was 2135143.5
$ws.Cells.Item(132, 11).Value = 2800.8  # K132: was 3107.625
$ws.Cells.Item(132, 12).Value = 6405336.600000001  # L132: was 6405430.5
$ws.Cells.Item(132, 13).Value = -270.8000000000002  # M132: was -577.625
$ws.Cells.Item(132, 14).Value = -6410396.600000001  # N132: was -6410490.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1160.2307  # H22: was 1164.5333
$ws.Cells.Item(22, 9).Value = 910.375  # I22: was 966.8
$ws.Cells.Item(22, 11).Value = 910.375  # K22: was 966.8
$ws.Cells.Item(22, 13).Value = -615.375  # M22: was -671.8

$ws.Cells.Item(27, 8).Value = 1160.2307  # H27: was 1164.5333
$ws.Cells.Item(27, 9).Value = 910.375  # I27: was 966.8
$ws.Cells.Item(27, 11).Value = 910.375  # K27: was 966.8
$ws.Cells.Item(27, 13).Value = -803.375  # M27: was -859.8

$ws.Cells.Item(40, 8).Value = 866.6667  # H40: was 950
$ws.Cells.Item(40, 9).Value = 866.6667  # I40: was 950
$ws.Cells.Item(40, 11).Value = 866.6667  # K40: was 950
$ws.Cells.Item(40, 13).Value = -730.6667  # M40: was -814

$ws.Cells.Item(55, 8).Value = 1232.5312  # H55: was 1271.3549
$ws.Cells.Item(55, 9).Value = 1037.4706  # I55: was 1099.25
$ws.Cells.Item(55, 10).Value = 1453.6  # J55: was 1454.9333
$ws.Cells.Item(55, 11).Value = 1037.4706  # K55: was 1099.25
$ws.Cells.Item(55, 12).Value = 1453.6  # L55: was 1454.9333
$ws.Cells.Item(55, 13).Value = -864.4706000000001  # M55: was -926.25
$ws.Cells.Item(55, 14).Value = -1799.6  # N55: was -1800.9333

$ws.Cells.Item(61, 8).Value = 3661.6  # H61: was 3603.9565
$ws.Cells.Item(61, 9).Value = 3179.2273  # I61: was 3261.524
$ws.Cells.Item(61, 10).Value = 7199  # J61: was 7199.5
$ws.Cells.Item(61, 11).Value = 3179.2273  # K61: was 3261.524
$ws.Cells.Item(61, 12).Value = 7199  # L61: was 7199.5
$ws.Cells.Item(61, 13).Value = -2977.2273  # M61: was -3059.524
$ws.Cells.Item(61, 14).Value = -7603  # N61: was -7603.5

$ws.Cells.Item(68, 8).Value = 16857  # H68: was 18999.834
$ws.Cells.Item(68, 10).Value = 3500  # J68: was 3333.3333
$ws.Cells.Item(68, 12).Value = 3500  # L68: was 3333.3333
$ws.Cells.Item(68, 14).Value = -4998  # N68: was -4831.3333

$ws.Cells.Item(71, 8).Value = 16857  # H71: was 18999.834
$ws.Cells.Item(71, 10).Value = 3500  # J71: was 3333.3333
$ws.Cells.Item(71, 12).Value = 17500  # L71: was 16666.6665
$ws.Cells.Item(71, 14).Value = -24988  # N71: was -24154.6665

$ws.Cells.Item(113, 8).Value = 3661.6  # H113: was 3603.9565
$ws.Cells.Item(113, 9).Value = 3179.2273  # I113: was 3261.524
$ws.Cells.Item(113, 10).Value = 7199  # J113: was 7199.5
$ws.Cells.Item(113, 11).Value = 3179.2273  # K113: was 3261.524
$ws.Cells.Item(113, 12).Value = 7199  # L113: was 7199.5
$ws.Cells.Item(113, 13).Value = -1009.2273  # M113: was -1091.524
$ws.Cells.Item(113, 14).Value = -11539  # N113: was -11539.5

$ws.Cells.Item(122, 8).Value = 6691.1875  # H122: was 7325.643
$ws.Cells.Item(122, 9).Value = 8027.2  # I122: was 9471.5
$ws.Cells.Item(122, 11).Value = 24081.6  # K122: was 28414.5
$ws.Cells.Item(122, 13).Value = -21631.6  # M122: was -25964.5

$ws.Cells.Item(132, 8).Value = 2415271.5  # H132: was 2415306
$ws.Cells.Item(132, 9).Value = 10220.071  # I132: was 10698.615
$ws.Cells.Item(132, 10).Value = 4659986.5  # J132: was 4369049.5
$ws.Cells.Item(132, 11).Value = 30660.213  # K132: was 32095.845
$ws.Cells.Item(132, 12).Value = 13979959.5  # L132: was 13107148.5
$ws.Cells.Item(132, 13).Value = -28130.213  # M132: was -29565.845
$ws.Cells.Item(132, 14).Value = -13985019.5  # N132: was -13112208.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 3373.6875  # H113: was 3374.375
$ws.Cells.Item(113, 9).Value = 4456.5  # I113: was 4456.9165
$ws.Cells.Item(113, 10).Value = 125.25  # J113: was 126.75
$ws.Cells.Item(113, 11).Value = 13369.5  # K113: was 13370.7495
$ws.Cells.Item(113, 12).Value = 375.75  # L113: was 380.25
$ws.Cells.Item(113, 13).Value = -11199.5  # M113: was -11200.7495
$ws.Cells.Item(113, 14).Value = -4715.75  # N113: was -4720.25

$ws.Cells.Item(122, 8).Value = 2023  # H122: was 2027.6875
$ws.Cells.Item(122, 9).Value = 1777.6  # I122: was 1807.5
$ws.Cells.Item(122, 10).Value = 2513.8  # J122: was 2394.6667
$ws.Cells.Item(122, 11).Value = 5332.799999999999  # K122: was 5422.5
$ws.Cells.Item(122, 12).Value = 7541.400000000001  # L122: was 7184.000100000001
$ws.Cells.Item(122, 13).Value = -2882.799999999999  # M122: was -2972.5
$ws.Cells.Item(122, 14).Value = -12441.4  # N122: was -12084.0001

$ws.Cells.Item(132, 8).Value = 726627.7  # H132: was 908038.8
$ws.Cells.Item(132, 9).Value = 1713.6364  # I132: was 1877.8889
$ws.Cells.Item(132, 10).Value = 2720141.2  # J132: was 3626521.8
$ws.Cells.Item(132, 11).Value = 5140.9092  # K132: was 5633.6667
$ws.Cells.Item(132, 12).Value = 8160423.600000001  # L132: was 10879565.4
$ws.Cells.Item(132, 13).Value = -2610.9092  # M132: was -3103.6667
$ws.Cells.Item(132, 14).Value = -8165483.600000001  # N132: was -10884625.4

$ws.Cells.Item(136, 8).Value = 580039.6  # H136: was 522080.34
$ws.Cells.Item(136, 9).Value = 2474.5  # I136: was 2184.8572
$ws.Cells.Item(136, 11).Value = 7423.5  # K136: was 6554.571599999999
$ws.Cells.Item(136, 13).Value = -4873.5  # M136: was -4004.571599999999
